# Apply "show circle and nodes inside" edit to the daily sheet workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 29 currently only has B29 = "3" (day label). Fill in the rest of the
# entry: description, begin time, finished time and hours worked.
$ws.Cells.Item(29, 3).Value = "Format node.json, Show nearest nodes inside circle."
$ws.Cells.Item(29, 4).Value = 0.95833333333333337
$ws.Cells.Item(29, 5).Value = 0.041666666666666664
$ws.Cells.Item(29, 6).Value = 2

# New row 30 starts the next day-label group with "4".
$ws.Cells.Item(30, 2).Value = "4"

# Keep the selection / active cell consistent with the saved state.
$ws.Range("C28").Select()
